$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 1823
$ws.Range("F8").Value = 1276
$ws.Range("F10").Value = 459
$ws.Range("F12").Value = 2669
$ws.Range("F13").Value = 361
$ws.Range("F14").Value = 867
$ws.Range("F16").Value = 572
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 55
$ws.Range("F19").Value = 1554
$ws.Range("F20").Value = 26
$ws.Range("F21").Value = 1223
$ws.Range("F22").Value = 169
$ws.Range("F23").Value = 593
$ws.Range("F26").Value = 1402
$ws.Range("F27").Value = 946
$ws.Range("F28").Value = 1309
$ws.Range("F29").Value = 205
$ws.Range("F30").Value = 1262
$ws.Range("F31").Value = 422
$ws.Range("F32").Value = 138
$ws.Range("F35").Value = 1823
$ws.Range("F36").Value = 460
$ws.Range("F40").Value = 2243
$ws.Range("F43").Value = 2741
$ws.Range("F46").Value = 143

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 30
$ws.Range("F13").Value = 106527
$ws.Range("F17").Value = 61
$ws.Range("F18").Value = 61
$ws.Range("F22").Value = 270
$ws.Range("F26").Value = 58
$ws.Range("F27").Value = 58
$ws.Range("F30").Value = 215
$ws.Range("F32").Value = 41

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F10").Value = 914
$ws.Range("F12").Value = 585
$ws.Range("F15").Value = 1164

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 914
$ws.Range("F9").Value = 585
$ws.Range("F12").Value = 1823
$ws.Range("F14").Value = 1276
$ws.Range("F16").Value = 459
$ws.Range("F17").Value = 1164
$ws.Range("F18").Value = 2669
$ws.Range("F19").Value = 30
$ws.Range("F20").Value = 361
$ws.Range("F21").Value = 867
$ws.Range("F23").Value = 572
$ws.Range("F24").Value = 1554
$ws.Range("F27").Value = 1223
$ws.Range("F28").Value = 169
$ws.Range("F29").Value = 593
$ws.Range("F30").Value = 1402
$ws.Range("F32").Value = 1309
$ws.Range("F33").Value = 205
$ws.Range("F35").Value = 61
$ws.Range("F36").Value = 1262
$ws.Range("F37").Value = 422
$ws.Range("F40").Value = 1823
$ws.Range("F41").Value = 58
$ws.Range("F44").Value = 2243
$ws.Range("F47").Value = 2741
$ws.Range("F49").Value = 143
